$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.346.02"
$ws.Range("E2").Value = "  -1.72%  "

$ws.Range("D3").Value = "3.063.77"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'588.48"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "'154.94"
$ws.Range("E6").Value = "  +4.22%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "'0.536"
$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").Value = "3.060.85"
$ws.Range("E9").Value = "  -3.00%  "

$ws.Range("E10").Value = "  -4.05%  "

$ws.Range("D11").Value = "'5.81"
$ws.Range("E11").Value = "  -1.78%  "

$ws.Range("E12").Value = "  -2.53%  "

$ws.Range("D13").Value = "'36.82"
$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("E14").Value = "  -4.65%  "

$ws.Range("E15").Value = "  -2.19%  "

$ws.Range("D16").Value = "3.571.98"
$ws.Range("E16").Value = "  -3.06%  "

$ws.Range("D17").Value = "63.425.13"
$ws.Range("E17").Value = "  -1.22%  "

$ws.Range("D18").Value = "'7.11"
$ws.Range("E18").Value = "  -2.67%  "

$ws.Range("D19").Value = "3.064.67"
$ws.Range("E19").Value = "  -2.84%  "

$ws.Range("D20").Value = "'471.12"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("D21").Value = "'14.27"
$ws.Range("E21").Value = "  -2.03%  "

$ws.Range("D22").Value = "'0.703"
$ws.Range("E22").Value = "  -4.82%  "

$ws.Range("D23").Value = "'7.49"
$ws.Range("E23").Value = "  -2.34%  "

$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = "  -1.62%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'80.51"
$ws.Range("E25").Value = "  -1.31%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "'12.79"
$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("D27").Value = "'10.34"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("D29").Value = "'7.38"
$ws.Range("E29").Value = "  +1.22%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("E31").Value = "  -2.87%  "

$ws.Range("D32").Value = "'2.14"
$ws.Range("E32").Value = "  -4.97%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "'27.06"
$ws.Range("E33").Value = "  -3.45%  "

$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.111"
$ws.Range("E34").Value = "  -3.86%  "

$ws.Range("D35").Value = "0.0₃0816"
$ws.Range("E35").Value = "  -5.39%  "

$ws.Range("E36").Value = "  -2.51%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'5.97"
$ws.Range("E37").Value = "  -4.11%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  -1.03%  "

$ws.Range("E39").Value = "  -5.33%  "

$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'9.23"
$ws.Range("E40").Value = "  -1.17%  "

$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'50.43"
$ws.Range("E41").Value = "  -1.97%  "

$ws.Range("D42").Value = "'435.22"
$ws.Range("E42").Value = "  -6.96%  "

$ws.Range("D43").Value = "'0.284"
$ws.Range("E43").Value = "  -3.84%  "

$ws.Range("D44").Value = "'40.64"
$ws.Range("E44").Value = "  +1.97%  "

$ws.Range("E45").Value = "  +3.22%  "

$ws.Range("D46").Value = "'0.0358"
$ws.Range("E46").Value = "  -4.87%  "

$ws.Range("D47").Value = "2.792.51"
$ws.Range("E47").Value = "  -4.04%  "

$ws.Range("D48").Value = "'129.61"
$ws.Range("E48").Value = "  -2.67%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "'25.02"
$ws.Range("E50").Value = "  +2.59%  "

$ws.Range("D51").Value = "'2.21"
$ws.Range("E51").Value = "  -1.80%  "
